$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, matching style of existing headers (A1:E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy style from an existing header cell so new headers match formatting
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# New boolean values for row 2
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false
